# Update "想去人数" (F column) figures on both the "展览" and "全部类型"
# worksheets to reflect the latest generated data.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F3").Value = 1255
    $ws.Range("F4").Value = 1522
    $ws.Range("F5").Value = 59
    $ws.Range("F6").Value = 6152
}
